# Werkbestand Projectadministratie - add Tier1 summary sheet (from "Sluiten")
# and close out the old per-project "Finale check" / "Eindacties" columns.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Rename the "Sluiten" working sheet to "Tier1".
$ws2.Name = "Tier1"

# Drop the per-project "Actiepunten Bram" / "Eindacties" / "Finale check"
# columns from the Tier1 summary sheet - Tier1 only needs the first four
# columns now.
$ws2.Columns("E:G").Delete()

# Sheet1 keeps its selection but is no longer the active tab - move it to I5
# (matches where the cursor is left after working the sheet).
$ws1.Range("I5").Select() | Out-Null

# Tier1 becomes the active sheet/tab, cursor parked on D1.
$ws2.Range("D1").Select() | Out-Null

Write-Output "applied edits"
